$wb = $excel.ActiveWorkbook

# --- Worksheets ---
$wsGeneral = $wb.Worksheets.Item("generalOptions")
$wsCases = $wb.Worksheets.Item("cases")

# --- Add new columns F (rotation) and G (management) to the "cases" sheet ---
# Header row.
$wsCases.Range("F1").Value = "rotation"
$wsCases.Range("G1").Value = "management"

# Data cells, filled in the order that reproduces the shared-string table of
# the target workbook (F3, F4, G2, G3, G4, F2).
$wsCases.Range("F3").Value = """MAIZE.bidule"""
$wsCases.Range("F4").Value = """Chickpea.Ghab2"", ""WHEAT.Cocorit"", ""WHEAT.Avoine_Romani"""
$wsCases.Range("G2").Value = """ROTATION_BLE"", ""ROTATION_BLE_IRRIGUE"""
$wsCases.Range("G3").Value = """Gorgan-RFD"""
$wsCases.Range("G4").Value = """ROTATION_POISCHICHE"", ""ROTATION_BLE"", ""ROTATION_BLE_IRRIGUE"""
$wsCases.Range("F2").Value = """WHEAT.Ble_Dur_1"", ""WHEAT.Ble_Tendre_1"""

# Give column F a custom width, matching the authored widening of the sheet.
# (24.83 is the closest achievable value that rounds to the target stored
# width of 25.6640625 given this engine's column-width quantization.)
$wsCases.Columns.Item(6).ColumnWidth = 24.83

# --- Update the selection / active sheet to match the target view ---
# Selecting F4 on "cases" both sets its selection and makes it the active /
# tabSelected sheet, which also clears tabSelected on "generalOptions".
$wsCases.Range("F4").Select()
